$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "29.417.38"
    "E2" = "  +0.01%  "
    "D3" = "1.851.77"
    "E3" = "  +0.21%  "
    "E4" = "  +0.10%  "
    "D5" = "241.02"
    "E5" = "  +0.18%  "
    "D6" = "0.6299"
    "E6" = "  -0.08%  "
    "D8" = "0.07686"
    "E8" = "  +1.63%  "
    "D9" = "0.2939"
    "E9" = "  -0.51%  "
    "D10" = "24.57"
    "E10" = "  +0.51%  "
    "D11" = "0.07750"
    "E11" = "  +0.75%  "
    "D12" = "1.854.49"
    "E12" = "  +0.41%  "
    "D13" = "5.031"
    "E13" = "  +0.77%  "
    "D14" = "0.00001094"
    "E14" = "  +8.93%  "
    "D15" = "0.6814"
    "E15" = "  -0.48%  "
    "D16" = "83.71"
    "E16" = "  +0.77%  "
    "D17" = "2.098.86"
    "E17" = "  -0.18%  "
    "E18" = "  +0.35%  "
    "D19" = "29.459.54"
    "E19" = "  +0.08%  "
    "D20" = "229.33"
    "E20" = "  +0.66%  "
    "D21" = "12.49"
    "E21" = "  +0.10%  "
    "D22" = "1.001"
    "E22" = "  +0.10%  "
    "E23" = "  -1.25%  "
    "D24" = "1.000"
    "E24" = "  +0.04%  "
    "D25" = "157.16"
    "E25" = "  +0.10%  "
    "D26" = "0.1390"
    "E26" = "  -0.14%  "
    "D27" = "8.393"
    "E27" = "  +0.23%  "
    "D28" = "17.68"
    "E28" = "  +0.04%  "
    "D29" = "1.313"
    "E29" = "  +3.70%  "
    "E30" = "  -0.24%  "
    "D31" = "0.05718"
    "E31" = "  +0.27%  "
    "D32" = "4.132"
    "E32" = "  +0.24%  "
    "D33" = "4.053"
    "E33" = "  +0.78%  "
    "E34" = "  +0.40%  "
    "D35" = "1.163"
    "E35" = "  +0.73%  "
    "D36" = "0.7090"
    "E36" = "  -0.53%  "
    "D37" = "2.588"
    "E37" = "  -0.05%  "
    "D38" = "2.779"
    "E38" = "  +0.04%  "
    "D39" = "0.01791"
    "E39" = "  -0.94%  "
    "D40" = "1.219.19"
    "E40" = "  -2.54%  "
    "D41" = "6.479"
    "E41" = "  +4.91%  "
    "D42" = "0.9080"
    "E42" = "  +0.23%  "
    "D43" = "1.001"
    "E43" = "  +0.05%  "
    "B44" = "Quant"
    "C44" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "D44" = "101.62"
    "E44" = "  +0.23%  "
    "B45" = "Aave"
    "C45" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D45" = "66.34"
    "E45" = "  +0.28%  "
    "B46" = "BabyDogeCoin"
    "C46" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D46" = "0.00000000120"
    "E46" = "  +1.91%  "
    "B47" = "Aptos"
    "C47" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D47" = "7.130"
    "E47" = "  -0.01%  "
    "B48" = "TheSandbox"
    "C48" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D48" = "0.4020"
    "E48" = "  +0.10%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D49" = "9.017"
    "E49" = "  -0.65%  "
    "D50" = "1.684"
    "E50" = "  +0.17%  "
    "B51" = "Algorand"
    "C51" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D51" = "0.1131"
    "E51" = "  +0.99%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}